# edit.ps1 - apply the commit's text edits to FraCVWalidRashad.docx
#
# Summary of changes (from the unified diff):
#  1. "Apprentissage rapide, ..." -> "Apprentissage extrêmement rapide, ..."
#     (new text split across 3 runs)
#  2. "Connaissance de nombreux langages informatiques (HTML, CSS, Javascript,
#     VBA, Python, etc.) " -> "Expertise en de nombreux langages informatiques
#     (HTML, CSS, Javascript, C, VBA, Matlab, Python, Assembleur, ...) "
#     (new text split across 8 runs)
#  3. "Université d'Ottawa, Ottawa (Ontario) " - collapse the proofErr-split
#     runs into a single run (text unchanged)
#  4. "Bourse d'études françaises " - same kind of run collapse
#  5. "Garde côtière canadienne " - same kind of run collapse
#  6. "Stagiaire en Compatibilité Électromagnétique" - same kind of run collapse
#  7. Add <w:lang w:val="fr-CA"/> to the last (empty) paragraph's mark
#     run-properties.

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $containsText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($containsText)) {
            return $i
        }
    }
    return -1
}

# Replace the whole text of the paragraph $paraIndex (previously found via
# Find-ParaIndex) that matches $oldText with $newText, as a single run.
function Replace-ParaText($doc, $paraIndex, $oldText, $newText) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.ClearFormatting()
    $r.Find.Text = $oldText
    $r.Find.Forward = $true
    $r.Find.MatchCase = $true
    $r.Find.MatchWildcards = $false
    $found = $r.Find.Execute()
    if (-not $found) {
        Write-Output ("NOT FOUND (replace): " + $oldText)
        return
    }
    $r.Text = $newText
}

# Force a run boundary around $text inside paragraph $paraIndex by toggling
# Bold on and back off (no visible formatting change, but it splits the run
# so the saved OOXML has $text as its own <w:r>).
function Split-AtTextInPara($doc, $paraIndex, $text) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.ClearFormatting()
    $r.Find.Text = $text
    $r.Find.Forward = $true
    $r.Find.MatchCase = $true
    $r.Find.MatchWildcards = $false
    $found = $r.Find.Execute()
    if (-not $found) {
        Write-Output ("NOT FOUND (split): " + $text)
        return
    }
    $r.Bold = 1
    $r.Bold = 0
}

# --- 1. "Apprentissage rapide, ..." ---------------------------------------
$idx = Find-ParaIndex $d "Apprentissage rapide"
Replace-ParaText $d $idx `
    "Apprentissage rapide, excellentes compétences en communication" `
    "Apprentissage extrêmement rapide, excellentes compétences en communication"
Split-AtTextInPara $d $idx "extrêmement "

# --- 2. "Connaissance de nombreux langages informatiques ..." -------------
$idx = Find-ParaIndex $d "Connaissance de nombreux"
Replace-ParaText $d $idx `
    "Connaissance de nombreux langages informatiques (HTML, CSS, Javascript, VBA, Python, etc.) " `
    "Expertise en de nombreux langages informatiques (HTML, CSS, Javascript, C, VBA, Matlab, Python, Assembleur, ...) "
Split-AtTextInPara $d $idx "nombreux langages informatiques (HTML, CSS, Javascript, "
Split-AtTextInPara $d $idx "C, "
Split-AtTextInPara $d $idx "VBA, "
Split-AtTextInPara $d $idx "Matlab, "
Split-AtTextInPara $d $idx "Python, "
Split-AtTextInPara $d $idx "Assembleur, .."
Split-AtTextInPara $d $idx ".) "

# --- 3. "Université d'Ottawa, Ottawa (Ontario) " ---------------------------
$idx = Find-ParaIndex $d "Universit"
Replace-ParaText $d $idx `
    "Université d'Ottawa, Ottawa (Ontario) " `
    "Université d'Ottawa, Ottawa (Ontario) X"
Replace-ParaText $d $idx `
    "Université d'Ottawa, Ottawa (Ontario) X" `
    "Université d'Ottawa, Ottawa (Ontario) "

# --- 4. "Bourse d'études françaises " --------------------------------------
$idx = Find-ParaIndex $d "Bourse"
Replace-ParaText $d $idx `
    "Bourse d'études françaises " `
    "Bourse d'études françaises X"
Replace-ParaText $d $idx `
    "Bourse d'études françaises X" `
    "Bourse d'études françaises "

# --- 5. "Garde côtière canadienne " ----------------------------------------
$idx = Find-ParaIndex $d "Garde"
Replace-ParaText $d $idx `
    "Garde côtière canadienne " `
    "Garde côtière canadienne X"
Replace-ParaText $d $idx `
    "Garde côtière canadienne X" `
    "Garde côtière canadienne "

# --- 6. "Stagiaire en Compatibilité Électromagnétique" --------------------
$idx = Find-ParaIndex $d "Stagiaire"
Replace-ParaText $d $idx `
    "Stagiaire en Compatibilité Électromagnétique" `
    "Stagiaire en Compatibilité Électromagnétique X"
Replace-ParaText $d $idx `
    "Stagiaire en Compatibilité Électromagnétique X" `
    "Stagiaire en Compatibilité Électromagnétique"

# --- 7. Add fr-CA language to the last (empty) paragraph mark -------------
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count - 1)
if ($last.Range.Text -ne "") {
    $last = $d.Paragraphs.Item($count)
}
$last.Range.LanguageID = "fr-CA"

Write-Output "done"
